$d = $word.ActiveDocument

# Step 1: fix the typo "Shoppe" -> "Shop" within the sentence.
$rng = $d.Content
$rng.Find.Execute("Book Shoppe,", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Book Shop,", 2) | Out-Null

# Step 2: split the run right after "...For online Book Shop" by toggling a
# formatting property on and back off. Word (and this runtime) breaks the
# paragraph's single run into two runs at the boundaries of a formatting
# change; since we immediately revert the property, both resulting runs end
# up with identical (original) run properties, matching the target markup.
$split = $d.Content
$split.Find.Execute("For online Book Shop", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$split.Bold = 1
$split.Bold = 0
